# "Generate Report for Handoff"
#
# The localization-status workbook tracks, per source file and per target
# locale, the timestamp of the most recent handoff. This run generates a
# fresh handoff report for the file
#   190b097a-fef3-43b3-a592-30c81669d904.md
# (row 4 on every sheet), bumping its "Latest Handoff Date/Datetime" to the
# new report's timestamps while leaving every other row/column untouched.

$wb = $excel.ActiveWorkbook

# Overview sheet: column D is "Latest Handoff Date" (file-level rollup).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("D4").Value = "2016-03-21 04:08:54"

# zh-cn sheet: column E is "Latest Handoff Datetime" for the zh-cn handoff.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "2016-03-21 04:08:45"

# de-de sheet: column E is "Latest Handoff Datetime" for the de-de handoff.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "2016-03-21 04:08:54"
